$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates
$ws.Range("B2").Value = 36.567242114214963
$ws.Range("C2").Value = 7.8053612560655665
$ws.Range("D2").Value = 10.971728364470863
$ws.Range("E2").Value = 2.1180024452172574

# Row 3 updates
$ws.Range("B3").Value = 53.193297954930813
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = -11.490467067090265
$ws.Range("E3").Value = 11.223432451595187

# Update selection to reflect new active range B1:E3
$ws.Range("B1:E3").Select()
